# Insert a new weekly price record for Espinaca (Vega Modelo de Temuco)
# as row 279, shifting the existing historical rows (old 279-309) down by one
# (new 280-310). This mirrors a new data point being prepended to the series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(279).Insert()

$ws.Range("A279").Value = 10
$ws.Range("B279").Value = "Vega Modelo de Temuco"
$ws.Range("C279").Value = "La Araucanía"
$ws.Range("D279").Value = 45194
$ws.Range("E279").Value = 9
$ws.Range("F279").Value = 100112012
$ws.Range("G279").Value = "Espinaca"
$ws.Range("H279").Value = "Sin especificar"
$ws.Range("I279").Value = "Primera"
$ws.Range("J279").Value = 50
$ws.Range("K279").Value = 10000
$ws.Range("L279").Value = 10000
$ws.Range("M279").Value = 10000
$ws.Range("N279").Value = "`$/docena de atados"
$ws.Range("O279").Value = "Región de La Araucanía"
$ws.Range("P279").Value = 3333
$ws.Range("Q279").Value = 3
$ws.Range("R279").Value = "Hortaliza"
